$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pure numeric refresh (countries whose row/rank did not change) ---
$ws.Range("B4").Value  = 1668493
$ws.Range("C4").Value  = 1665
$ws.Range("D4").Value  = 446982
$ws.Range("E4").Value  = 1122805
$ws.Range("G4").Value  = 23
$ws.Range("H4").Value  = 98706

$ws.Range("B14").Value = 133725
$ws.Range("C14").Value = 2302
$ws.Range("E14").Value = 74951
$ws.Range("G14").Value = 41
$ws.Range("H14").Value = 3909

$ws.Range("B55").Value = 8349
$ws.Range("C55").Value = 3
$ws.Range("E55").Value = 387

$ws.Range("D154").Value = 122
$ws.Range("E154").Value = 73

# --- Rows 89/90: Gabon/Cuba swap rank (Cuba now ahead of Gabon) ---
$ws.Range("A89").Value = "Cuba"
$ws.Range("B89").Value = 1941
$ws.Range("C89").Value = 10
$ws.Range("D89").Value = 1689
$ws.Range("E89").Value = 170
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 82

$ws.Range("A90").Value = "Gabon"
$ws.Range("B90").Value = 1934
$ws.Range("D90").Value = 459
$ws.Range("E90").Value = 1463
$ws.Range("H90").Value = 12

# --- Rows 198/199/200: Belice, Santa Lucia, Nueva Caledonia re-ranked ---
$ws.Range("A198").Value = "Nueva Caledonia"
$ws.Range("D198").Value = 18
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Belice"
$ws.Range("D199").Value = 16
$ws.Range("H199").Value = 2

$ws.Range("A200").Value = "Santa Lucia"

# --- Rows 209/210/211: Montserrat, Seychelles, Groenlandia re-ranked ---
$ws.Range("A209").Value = "Seychelles"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = "Groenlandia"

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# --- Rows 214/216: Sahara Occidental <-> Bonaire, San Eustaquio y Saba ---
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "Sahara Occidental"

# --- Timestamp caption ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 16:05"
